# Update column G ("K" = strikeouts) values to reflect regenerated save_data
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 2
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 1
    18 = 2
    20 = 1
    21 = 0
    22 = 0
    23 = 0
    25 = 3
    26 = 2
    27 = 2
    29 = 1
    30 = 1
    31 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
